$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 174, pushing the former rows 175-180
# down to 177-182. The values of the shifted rows stay intact automatically.
$ws.Range("A175:A176").EntireRow.Insert()

# --- New row 175: Red Globe, Provincia del Elquí ---
$ws.Range("A175").Value = 4
$ws.Range("B175").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C175").Value = "Los Lagos"
$ws.Range("D175").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D175").Value = 44568
$ws.Range("E175").Value = 10
$ws.Range("F175").Value = "Fruta"
$ws.Range("G175").Value = 100109
$ws.Range("H175").Value = "Uva"
$ws.Range("I175").Value = 100109001
$ws.Range("J175").Value = "Uva"
$ws.Range("K175").Value = "Red Globe"
$ws.Range("L175").Value = "Primera"
$ws.Range("M175").Value = 300
$ws.Range("N175").Value = 17000
$ws.Range("O175").Value = 18000
$ws.Range("P175").Value = 17500
$ws.Range("Q175").Value = "$/bandeja 8 kilos"
$ws.Range("R175").Value = "Provincia del Elquí"
$ws.Range("S175").Value = 2188
$ws.Range("T175").Value = 8

# --- New row 176: Superior Seedless, Provincia de Limarí ---
$ws.Range("A176").Value = 4
$ws.Range("B176").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C176").Value = "Los Lagos"
$ws.Range("D176").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D176").Value = 44568
$ws.Range("E176").Value = 10
$ws.Range("F176").Value = "Fruta"
$ws.Range("G176").Value = 100109
$ws.Range("H176").Value = "Uva"
$ws.Range("I176").Value = 100109001
$ws.Range("J176").Value = "Uva"
$ws.Range("K176").Value = "Superior Seedless"
$ws.Range("L176").Value = "Primera"
$ws.Range("M176").Value = 300
$ws.Range("N176").Value = 15000
$ws.Range("O176").Value = 16000
$ws.Range("P176").Value = 15500
$ws.Range("Q176").Value = "$/bandeja 8 kilos"
$ws.Range("R176").Value = "Provincia de Limarí"
$ws.Range("S176").Value = 1938
$ws.Range("T176").Value = 8
